$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-5 from 2023-10-25 (45224) to 2023-11-03 (45233)
$ws.Range("C2:C5").Value = 45233
